$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.368.12'
$ws.Range("E2").Value = '  +2.09%  '
$ws.Range("D3").Value = '3.420.36'
$ws.Range("E3").Value = '  +0.89%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '566.84'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.43%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '178.28'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.71%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.633'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.34%  '
$ws.Range("D8").Value = '3.413.21'
$ws.Range("E8").Value = '  +0.85%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.03%  '
$ws.Range("E10").Value = '  +4.66%  '
$ws.Range("E11").Value = '  +1.21%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.53'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.01%  '
$ws.Range("E13").Value = '  +0.69%  '
$ws.Range("E14").Value = '  +2.27%  '
$ws.Range("D15").Value = '3.947.49'
$ws.Range("E15").Value = '  +0.48%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '18.33'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.16%  '
$ws.Range("E17").Value = '  +0.93%  '
$ws.Range("D18").Value = '3.414.57'
$ws.Range("E18").Value = '  +1.06%  '
$ws.Range("D19").Value = '66.195.18'
$ws.Range("E19").Value = '  +2.04%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.98'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.31%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.01'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.32%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '466.91'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.96'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.41%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '14.79'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +9.12%  '
$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.14'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.32%  '
$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '90.01'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.38%  '
$ws.Range("E27").Value = '  +1.08%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.78'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.00%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.84'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.73%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '31.45'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.42%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.72'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.33%  '
$ws.Range("E32").Value = '  +0.38%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '583.26'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.78%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '62.54'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.75%  '
$ws.Range("E35").Value = '  +0.82%  '
$ws.Range("E36").Value = '  -0.04%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.144'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.01%  '
$ws.Range("E38").Value = '  +0.35%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.60'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.65%  '
$ws.Range("E40").Value = '  +3.27%  '
$ws.Range("D41").Value = '0.0₃0764'
$ws.Range("E41").Value = '  +2.01%  '
$ws.Range("D42").Value = '3.131.95'
$ws.Range("E42").Value = '  +1.16%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.89'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.47%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0423'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.88%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.50'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.45%  '
$ws.Range("E46").Value = '  -0.31%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.17'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.48%  '
$ws.Range("B48").Value = 'dogwifhat'
$ws.Range("C48").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.64'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +13.83%  '
$ws.Range("B49").Value = 'FirstDigitalUSD'
$ws.Range("C49").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.00'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.10%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '141.38'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.39%  '
$ws.Range("B51").Value = 'THORChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.58'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.95%  '
